$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1 (Word-index): 16÷3=, 77÷3=, 32÷3=, 81÷3=, 51÷2=
#                  -> 58÷8=, 67÷6=, 39÷8=, 10÷4=, 27÷8=
$row = $t.Rows.Item(1)
$row.Cells.Item(1).Range.Text = "58÷8="
$row.Cells.Item(2).Range.Text = "67÷6="
$row.Cells.Item(3).Range.Text = "39÷8="
$row.Cells.Item(4).Range.Text = "10÷4="
$row.Cells.Item(5).Range.Text = "27÷8="

# Row 5 (Word-index): 17÷7=, 69÷8=, 85÷5=, 95÷2=, 98÷3=
#   -> cell1 text changes; a new cell (94÷3=) is inserted right after cell1;
#      old cell2 (69÷8=) shifts to position 3 unchanged; old cell3 (85÷5=) becomes
#      position 4 with new text 91÷6=; old cell4 (95÷2=) becomes position 5 with new
#      text 68÷6=; old cell5 (98÷3=) is removed.
#   Net result (5 cells total, unchanged count): 29÷5=, 94÷3=, 69÷8=, 91÷6=, 68÷6=
$row = $t.Rows.Item(5)
$row.Cells.Item(1).Range.Text = "29÷5="
$row.Cells.Item(2).Range.Text = "94÷3="
$row.Cells.Item(3).Range.Text = "69÷8="
$row.Cells.Item(4).Range.Text = "91÷6="
$row.Cells.Item(5).Range.Text = "68÷6="

# Row 9 (Word-index): 11÷6=, 95÷2=, 72÷5=, 61÷4=, 98÷2=
#                  -> 70÷4=, 83÷8=, 57÷7=, 49÷6=, 27÷3=
$row = $t.Rows.Item(9)
$row.Cells.Item(1).Range.Text = "70÷4="
$row.Cells.Item(2).Range.Text = "83÷8="
$row.Cells.Item(3).Range.Text = "57÷7="
$row.Cells.Item(4).Range.Text = "49÷6="
$row.Cells.Item(5).Range.Text = "27÷3="

# Row 13 (Word-index): 22÷4=, 12÷5=, 16÷7=, 46÷9=, 59÷9=
#                   -> 88÷3=, 35÷2=, 75÷4=, 43÷4=, 81÷7=
$row = $t.Rows.Item(13)
$row.Cells.Item(1).Range.Text = "88÷3="
$row.Cells.Item(2).Range.Text = "35÷2="
$row.Cells.Item(3).Range.Text = "75÷4="
$row.Cells.Item(4).Range.Text = "43÷4="
$row.Cells.Item(5).Range.Text = "81÷7="

# Row 17 (Word-index): 39÷6=, 34÷8=, 19÷7=, 64÷9=, 69÷8=
#                   -> 15÷4=, 69÷2=, 55÷9=, 23÷2=, 52÷9=
$row = $t.Rows.Item(17)
$row.Cells.Item(1).Range.Text = "15÷4="
$row.Cells.Item(2).Range.Text = "69÷2="
$row.Cells.Item(3).Range.Text = "55÷9="
$row.Cells.Item(4).Range.Text = "23÷2="
$row.Cells.Item(5).Range.Text = "52÷9="
